$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Root-cause content edit: the "Character stats" window's size value changes
# from 0e20 to 0e1c (leading apostrophe keeps it quote-prefixed text, matching
# the cell's existing scientific-notation display style instead of Excel
# assigning a brand new style record for it).
$ws.Range("E17").Value2 = "'0e1c"

# Downstream re-pointing of two PatchWindow start-address formulas so the
# "Hapsby travel" and "Use, Equip, Drop" windows chain off their new
# neighbours.
$ws.Range("B13").Formula = "=B12"
$ws.Range("B18").Formula = "=C20"

# Restore the author's on-screen selection.
[void]$ws.Range("K17").Select()
